$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.7
$ws.Range("D3").Value = -7.7
$ws.Range("D5").Value = -7.953999999999999
$ws.Range("A9").Value = -20.912
$ws.Range("D11").Value = -8.15
$ws.Range("D12").Value = -8.088999999999999
$ws.Range("A13").Value = -22.005
$ws.Range("A16").Value = -20.86
$ws.Range("A18").Value = -21.868
$ws.Range("A20").Value = -21.664
$ws.Range("D21").Value = -7.951000000000001
